$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 5942649
$ws.Range("C4").Value = 27019
$ws.Range("D4").Value = 3232526
$ws.Range("E4").Value = 2528137
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 872
$ws.Range("H4").Value = 181986

# Row 8
$ws.Range("B8").Value = 613017
$ws.Range("C8").Value = 1567
$ws.Range("D8").Value = 520381
$ws.Range("E8").Value = 79328
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 149
$ws.Range("H8").Value = 13308

# Row 23
$ws.Range("B23").Value = 237568
$ws.Range("C23").Value = 1451
$ws.Range("D23").Value = 209600
$ws.Range("E23").Value = 18623
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 9345

# Row 27
$ws.Range("B27").Value = 125834
$ws.Range("C27").Value = 187
$ws.Range("D27").Value = 111890
$ws.Range("E27").Value = 4855
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9089

# Row 32
$ws.Range("B32").Value = 106460
$ws.Range("C32").Value = 1988
$ws.Range("D32").Value = 83822
$ws.Range("E32").Value = 21779
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 12
$ws.Range("H32").Value = 859

# Row 65
$ws.Range("B65").Value = 35305
$ws.Range("C65").Value = 842
$ws.Range("D65").Value = 13072
$ws.Range("E65").Value = 21857
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 14
$ws.Range("H65").Value = 376

# Row 103
$ws.Range("A103").Value = "Maldivas"
$ws.Range("B103").Value = 7047
$ws.Range("C103").Value = 135
$ws.Range("D103").Value = 4439
$ws.Range("E103").Value = 2580
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 28

# Row 104
$ws.Range("A104").Value = "Mauritania"
$ws.Range("B104").Value = 6928
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 6282
$ws.Range("E104").Value = 488
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 158

# Row 117
$ws.Range("B117").Value = 3744
$ws.Range("C117").Value = 27
$ws.Range("D117").Value = 3102
$ws.Range("E117").Value = 551
$ws.Range("F117").Value = 0

# Row 126
$ws.Range("A126").Value = "Tunez"
$ws.Range("B126").Value = 3069
$ws.Range("C126").Value = 176
$ws.Range("D126").Value = 1456
$ws.Range("E126").Value = 1542
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 71

# Row 127
$ws.Range("A127").Value = "Sri Lanka"
$ws.Range("B127").Value = 2971
$ws.Range("C127").Value = 12
$ws.Range("D127").Value = 2816
$ws.Range("E127").Value = 143
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 12

# Row 146
$ws.Range("B146").Value = 1670
$ws.Range("C146").Value = 42
$ws.Range("D146").Value = 565
$ws.Range("E146").Value = 1097
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 8

# Row 156
$ws.Range("B156").Value = 1292
$ws.Range("C156").Value = 48
$ws.Range("D156").Value = 692
$ws.Range("E156").Value = 594
$ws.Range("F156").Value = 0

# Row 166
$ws.Range("B166").Value = 730
$ws.Range("C166").Value = 17
$ws.Range("D166").Value = 48
$ws.Range("E166").Value = 672
$ws.Range("F166").Value = 0

# Row 192
$ws.Range("B192").Value = 136
$ws.Range("C192").Value = 4
$ws.Range("D192").Value = 127
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = 0

# Row 203
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 26
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 25
$ws.Range("E203").Value = 1
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("B204").Value = 26
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 25
$ws.Range("E204").Value = 1
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0
